$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.317.23"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "1.602.65"
$ws.Range("E3").Value = "  +1.17%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("E9").Value = "  -0.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("E11").Value = "  +0.82%  "
$ws.Range("D12").Value = "1.828.89"
$ws.Range("E12").Value = "  +1.19%  "
$ws.Range("D13").Value = "1.604.13"
$ws.Range("E13").Value = "  +1.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.506"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.66"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("D17").Value = "26.315.27"
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "227.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.85%  "
$ws.Range("D19").Value = "0.0₃0722"
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.90%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("E22").Value = "  +1.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.26%  "
$ws.Range("E28").Value = "  +1.57%  "
$ws.Range("E29").Value = "  +2.30%  "
$ws.Range("E30").Value = "  -0.38%  "
$ws.Range("E31").Value = "  +0.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.79%  "
$ws.Range("D33").Value = "1.444.74"
$ws.Range("E33").Value = "  +7.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.57%  "
$ws.Range("E36").Value = "  +1.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.565"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.38%  "
$ws.Range("E38").Value = "  -0.73%  "
$ws.Range("E39").Value = "  +1.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.82"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.60%  "
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("E42").Value = "  +2.36%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.922"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.71%  "
$ws.Range("D44").Value = "1.740.83"
$ws.Range("E44").Value = "  +1.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.759"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.85"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.10%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.78%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0500"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₇0969"
$ws.Range("E50").Value = "  -5.45%  "
$ws.Range("E51").Value = "  -3.37%  "

Write-Host "Update complete"